# Rename the "Latest Payment Date" column header to "Payment Date"
# on the Expense sheet (cell E1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expense")

$ws.Range("E1").Value = "Payment Date"

# Reflect the cursor/selection ending up on E2 after the edit.
$ws.Activate()
$ws.Range("E2").Select()
